$d = $word.ActiveDocument

# Update the date paragraph
$d.Content.Find.Execute("2024-05-04 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-05 Sunday", 2)

# Update the table of equations, in row-major order (20 rows x 5 cols)
$newValues = @(
  "96-49=47",
  "61+22=83",
  "62-4=58",
  "14+74=88",
  "14+61=75",
  "8+20=28",
  "44+26=70",
  "24+66=90",
  "18+64=82",
  "24+17=41",
  "99-73=26",
  "59+5=64",
  "98-97=1",
  "48+7=55",
  "26-14=12",
  "73-23=50",
  "37-8=29",
  "70-55=15",
  "76-50=26",
  "69-54=15",
  "79-73=6",
  "9-4=5",
  "95-18=77",
  "49-7=42",
  "95-42=53",
  "13+73=86",
  "33+16=49",
  "16+28=44",
  "17+9=26",
  "47-43=4",
  "90-44=46",
  "42+43=85",
  "32-2=30",
  "46+45=91",
  "4+43=47",
  "96-42=54",
  "47+25=72",
  "96-21=75",
  "12+8=20",
  "66-54=12",
  "76-48=28",
  "42-9=33",
  "17+56=73",
  "7+65=72",
  "67+20=87",
  "55+32=87",
  "78-6=72",
  "46+3=49",
  "24+75=99",
  "14+79=93",
  "43+39=82",
  "77-20=57",
  "18+31=49",
  "17+32=49",
  "27-7=20",
  "47+14=61",
  "79-66=13",
  "43+26=69",
  "86-24=62",
  "64-3=61",
  "46+25=71",
  "66+29=95",
  "6+33=39",
  "34+54=88",
  "77-69=8",
  "0+58=58",
  "84-68=16",
  "40+47=87",
  "91-79=12",
  "9+22=31",
  "92+4=96",
  "54-29=25",
  "96-26=70",
  "56-28=28",
  "54+23=77",
  "0+29=29",
  "2+91=93",
  "5+15=20",
  "28-18=10",
  "20+47=67",
  "21+49=70",
  "43+49=92",
  "69+1=70",
  "45+30=75",
  "48+12=60",
  "70+26=96",
  "93-41=52",
  "66-28=38",
  "75-10=65",
  "14+46=60",
  "0+75=75",
  "74-55=19",
  "80-55=25",
  "58+38=96",
  "10+75=85",
  "67-2=65",
  "56-52=4",
  "13+23=36",
  "26-26=0",
  "41-33=8"
)

$t = $d.Tables.Item(1)
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
  for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $newValues[$idx]
    $idx++
  }
}

Write-Host "Done. Updated $idx cells."
